$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 46

$newLink = "https://www.medpagetoday.com/meetingcoverage/ims/117570"
$newKeyword = "MGUS"
$newTitle = "Linvoseltamab Promising in High-Risk Smoldering Multiple Myeloma"

# Fill in the new data row
$ws.Cells.Item($newRow, 2).Value = $newKeyword
$ws.Cells.Item($newRow, 3).Value = $newTitle

# Add the hyperlink for the new row (A46); Excel will set the text and apply the Hyperlink style
$ws.Hyperlinks.Add($ws.Cells.Item($newRow, 1), $newLink, [Type]::Missing, [Type]::Missing, $newLink) | Out-Null
$ws.Cells.Item($newRow, 1).Style = "Hyperlink"
